# Actualización automática 2025-11-28 15:30:09
#
# Applies new sales figures recorded for several clients of
# "LINDAO ZUÑIGA BRYAN JOSE" across the three report sheets:
#   - VENTAS POR GRUPO      (sales broken down by product group)
#   - VENTA MENSUAL         (sales broken down by month, November column)
#   - CUMPLIMIENTO MENSUAL  (budget vs. actual compliance roll-up)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("K4").Value = 785.09
$wsGrupo.Range("M7").Value = 2637.26
$wsGrupo.Range("L16").Value = 537.34
$wsGrupo.Range("M16").Value = 8542.93
$wsGrupo.Range("M27").Value = 1891.36
$wsGrupo.Range("D35").Value = 732.6799999999999
$wsGrupo.Range("M35").Value = 1037.41
$wsGrupo.Range("M46").Value = 2179.48

# Summary row (row 60) - "N de 58" occurrence counters
$wsGrupo.Range("D60").Value = "3 de 58"
$wsGrupo.Range("K60").Value = "1 de 58"
$wsGrupo.Range("L60").Value = "4 de 58"
$wsGrupo.Range("M60").Value = "11 de 58"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 785.09
$wsMensual.Range("F7").Value = 2637.26
$wsMensual.Range("F16").Value = 9080.27
$wsMensual.Range("F27").Value = 1891.36
$wsMensual.Range("F35").Value = 1770.09
$wsMensual.Range("F46").Value = 2179.48
$wsMensual.Range("F60").Value = 53165.42

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 3471.56
$wsCumplimiento.Range("E3").Value = -1603.87
$wsCumplimiento.Range("F3").Value = 1.858745294990068

$wsCumplimiento.Range("D10").Value = 3377.81
$wsCumplimiento.Range("E10").Value = -2989.702016465608
$wsCumplimiento.Range("F10").Value = 8.703273684914231

$wsCumplimiento.Range("D11").Value = 14407.86
$wsCumplimiento.Range("E11").Value = -12966.94
$wsCumplimiento.Range("F11").Value = 9.999070038586458

$wsCumplimiento.Range("D12").Value = 36359.92
$wsCumplimiento.Range("E12").Value = 11681.08
$wsCumplimiento.Range("F12").Value = 0.756851855706584

$wsCumplimiento.Range("D14").Value = 61566.89999999999
$wsCumplimiento.Range("E14").Value = -3679.548035021517
$wsCumplimiento.Range("F14").Value = 1.063563937719031

Write-Output "Applied sales update for LINDAO ZUÑIGA BRYAN JOSE"
